# Fruta / hortaliza, semanal
# Insert 3 new weekly price records (rows) right before the current row 604,
# pushing the existing data down by 3 rows, and fill the new rows with the
# latest week's data (same product/location category, new date + prices).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 604, shifting rows 604:644 down to 607:647.
$ws.Rows("604:606").Insert()

# Constant columns for this product/category block (Terminal Hortofrutícola
# Agro Chillán, Ñuble, Fruta, Berries, Frutilla, Sin especificar, weekly=7).
$colA = 7
$colB = "Terminal Hortofrutícola Agro Chillán"
$colC = "Ñuble"
$colE = 16
$colF = "Fruta"
$colG = 100101
$colH = "Berries"
$colI = 100112025
$colJ = "Frutilla"
$colK = "Sin especificar"
$colT = 7

function Set-RowData {
    param($row, $d, $l, $m, $n, $o, $p, $q, $r, $s)

    $ws.Range("A$row").Value = $colA
    $ws.Range("B$row").Value = $colB
    $ws.Range("C$row").Value = $colC
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $colE
    $ws.Range("F$row").Value = $colF
    $ws.Range("G$row").Value = $colG
    $ws.Range("H$row").Value = $colH
    $ws.Range("I$row").Value = $colI
    $ws.Range("J$row").Value = $colJ
    $ws.Range("K$row").Value = $colK
    $ws.Range("L$row").Value = $l
    $ws.Range("M$row").Value = $m
    $ws.Range("N$row").Value = $n
    $ws.Range("O$row").Value = $o
    $ws.Range("P$row").Value = $p
    $ws.Range("Q$row").Value = $q
    $ws.Range("R$row").Value = $r
    $ws.Range("S$row").Value = $s
    $ws.Range("T$row").Value = $colT
}

Set-RowData 604 45265 "Especial" 100 12000 12000 12000 "$/bandeja 7 kilos" "Provincia de Melipilla" 1714
Set-RowData 605 45265 "Primera"  150 10000 10000 10000 "$/bandeja 7 kilos" "Provincia de Melipilla" 1429
Set-RowData 606 45265 "Segunda"  100  8000  8000  8000 "$/bandeja 7 kilos" "Provincia de Melipilla" 1143
